$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(5)

# Move the "Private LAN" textbox slightly to the right (x offset 3459377 -> 3468902 EMU)
$shape.Left = 273.1419

# Remove the "Private" run, leaving just "LAN" (drops the leading space too)
$shape.TextFrame.TextRange.Text = "LAN"
